# Auto-generated by analysis of the OOXML diff.
# Updates currentAveragePrice*, LevePrice*, LeveProfit* columns (H:N)
# for specific leve rows across the ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR sheets,
# mirroring a scheduled market-data refresh.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 161.44444
$ws.Range("I2").Value = 109
$ws.Range("J2").Value = 345
$ws.Range("K2").Value = 109
$ws.Range("L2").Value = 345
$ws.Range("M2").Value = 4
$ws.Range("N2").Value = -571
$ws.Range("H15").Value = 1092.2759
$ws.Range("I15").Value = 1092.2759
$ws.Range("K15").Value = 3276.8277
$ws.Range("M15").Value = -3107.8277
$ws.Range("H92").Value = 522.44446
$ws.Range("I92").Value = 409.38095
$ws.Range("J92").Value = 918.1667
$ws.Range("K92").Value = 409.38095
$ws.Range("L92").Value = 918.1667
$ws.Range("M92").Value = 838.61905
$ws.Range("N92").Value = -3414.1667
$ws.Range("H106").Value = 35561.855
$ws.Range("I106").Value = 19328
$ws.Range("K106").Value = 19328
$ws.Range("M106").Value = -18697
$ws.Range("H111").Value = 1703.85
$ws.Range("I111").Value = 1564.5
$ws.Range("K111").Value = 4693.5
$ws.Range("M111").Value = -1626.5
$ws.Range("H113").Value = 5788.0415
$ws.Range("I113").Value = 4862.4
$ws.Range("J113").Value = 7330.778
$ws.Range("K113").Value = 4862.4
$ws.Range("L113").Value = 7330.778
$ws.Range("M113").Value = -1608.4
$ws.Range("N113").Value = -13838.778
$ws.Range("H132").Value = 1408.7333
$ws.Range("I132").Value = 1416.0714
$ws.Range("K132").Value = 4248.2142
$ws.Range("M132").Value = -1718.2142
$ws.Range("H134").Value = 149999
$ws.Range("J134").Value = 149999
$ws.Range("L134").Value = 149999
$ws.Range("N134").Value = -160139
$ws.Range("H135").Value = 3149.7856
$ws.Range("I135").Value = 3149.7856
$ws.Range("K135").Value = 28348.0704
$ws.Range("M135").Value = -25813.0704
$ws.Range("H137").Value = 18183.234
$ws.Range("I137").Value = 9912.799999999999
$ws.Range("J137").Value = 29998.143
$ws.Range("K137").Value = 29738.4
$ws.Range("L137").Value = 89994.429
$ws.Range("M137").Value = -27188.4
$ws.Range("N137").Value = -95094.429
$ws.Range("H138").Value = 2674.7856
$ws.Range("I138").Value = 2985.9
$ws.Range("J138").Value = 1897
$ws.Range("K138").Value = 8957.700000000001
$ws.Range("L138").Value = 5691
$ws.Range("M138").Value = -3817.700000000001
$ws.Range("N138").Value = -15971
$ws.Range("H139").Value = 0
$ws.Range("J139").Value = 0
$ws.Range("L139").ClearContents()
$ws.Range("N139").Value = 0
$ws.Range("H141").Value = 3814.8
$ws.Range("I141").Value = 4358.6665
$ws.Range("K141").Value = 13075.9995
$ws.Range("M141").Value = -7895.999500000002

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4485.524
$ws.Range("I32").Value = 4485.524
$ws.Range("J32").Value = 0
$ws.Range("K32").Value = 4485.524
$ws.Range("L32").Value = 0
$ws.Range("M32").ClearContents()
$ws.Range("N32").Value = -4198.524
$ws.Range("H45").Value = 7267.7144
$ws.Range("I45").Value = 9964.385
$ws.Range("K45").Value = 9964.385
$ws.Range("M45").Value = -9587.385
$ws.Range("H61").Value = 11155.322
$ws.Range("I61").Value = 9815.406999999999
$ws.Range("J61").Value = 20199.75
$ws.Range("K61").Value = 9815.406999999999
$ws.Range("L61").Value = 20199.75
$ws.Range("M61").Value = -9603.406999999999
$ws.Range("N61").Value = -20623.75
$ws.Range("H74").Value = 250000
$ws.Range("I74").Value = 250000
$ws.Range("K74").Value = 250000
$ws.Range("M74").Value = -249126
$ws.Range("H77").Value = 250000
$ws.Range("I77").Value = 250000
$ws.Range("K77").Value = 1250000
$ws.Range("M77").Value = -1245632
$ws.Range("H102").Value = 1183.68
$ws.Range("I102").Value = 1183.68
$ws.Range("K102").Value = 1183.68
$ws.Range("M102").Value = 438.3199999999999
$ws.Range("H122").Value = 1710.6757
$ws.Range("I122").Value = 1572.9032
$ws.Range("J122").Value = 2422.5
$ws.Range("K122").Value = 4718.7096
$ws.Range("L122").Value = 7267.5
$ws.Range("M122").Value = -2268.7096
$ws.Range("N122").Value = -12167.5
$ws.Range("H132").Value = 6120.8823
$ws.Range("I132").Value = 3338
$ws.Range("J132").Value = 12799.8
$ws.Range("K132").Value = 10014
$ws.Range("L132").Value = 38399.39999999999
$ws.Range("M132").Value = -7484
$ws.Range("N132").Value = -43459.39999999999
$ws.Range("H136").Value = 11155.322
$ws.Range("I136").Value = 9815.406999999999
$ws.Range("J136").Value = 20199.75
$ws.Range("K136").Value = 29446.221
$ws.Range("L136").Value = 60599.25
$ws.Range("M136").Value = -26896.221
$ws.Range("N136").Value = -65699.25

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 4600.1816
$ws.Range("I105").Value = 4276.615
$ws.Range("J105").Value = 5067.5557
$ws.Range("K105").Value = 4276.615
$ws.Range("L105").Value = 5067.5557
$ws.Range("M105").Value = -2529.615
$ws.Range("N105").Value = -8561.555700000001

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2170.6296
$ws.Range("I31").Value = 1281.7059
$ws.Range("J31").Value = 3681.8
$ws.Range("K31").Value = 1281.7059
$ws.Range("L31").Value = 3681.8
$ws.Range("M31").Value = -986.7058999999999
$ws.Range("N31").Value = -4271.8
$ws.Range("H34").Value = 2170.6296
$ws.Range("I34").Value = 1281.7059
$ws.Range("J34").Value = 3681.8
$ws.Range("K34").Value = 1281.7059
$ws.Range("L34").Value = 3681.8
$ws.Range("M34").Value = -1079.7059
$ws.Range("N34").Value = -4085.8
$ws.Range("H58").Value = 5965.1665
$ws.Range("J58").Value = 10107.1
$ws.Range("L58").Value = 10107.1
$ws.Range("N58").Value = -10513.1
$ws.Range("H136").Value = 5965.1665
$ws.Range("J136").Value = 10107.1
$ws.Range("L136").Value = 30321.3
$ws.Range("N136").Value = -35421.3

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H120").Value = 1999.5
$ws.Range("I120").Value = 1999.5
$ws.Range("K120").Value = 5998.5
$ws.Range("M120").Value = -1160.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H93").Value = 42993
$ws.Range("I93").Value = 0
$ws.Range("K93").Value = 0
$ws.Range("M93").ClearContents()
$ws.Range("H97").Value = 1102.25
$ws.Range("I97").Value = 1048.7142
$ws.Range("J97").Value = 1177.2
$ws.Range("K97").Value = 1048.7142
$ws.Range("L97").Value = 1177.2
$ws.Range("M97").Value = -552.7141999999999
$ws.Range("N97").Value = -2169.2
$ws.Range("H102").Value = 1488.3704
$ws.Range("I102").Value = 1519.6364
$ws.Range("J102").Value = 1350.8
$ws.Range("K102").Value = 1519.6364
$ws.Range("L102").Value = 1350.8
$ws.Range("M102").Value = 102.3635999999999
$ws.Range("N102").Value = -4594.8
$ws.Range("H113").Value = 264575.12
$ws.Range("I113").Value = 288085.84
$ws.Range("K113").Value = 288085.84
$ws.Range("M113").Value = -285915.84
$ws.Range("H122").Value = 1673.5834
$ws.Range("I122").Value = 1655.2
$ws.Range("J122").Value = 1765.5
$ws.Range("K122").Value = 4965.6
$ws.Range("L122").Value = 5296.5
$ws.Range("M122").Value = -2515.6
$ws.Range("N122").Value = -10196.5
$ws.Range("H132").Value = 69212.664
$ws.Range("I132").Value = 37319
$ws.Range("K132").Value = 111957
$ws.Range("M132").Value = -109427

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 9048.866
$ws.Range("I7").Value = 7374.5
$ws.Range("J7").Value = 12397.6
$ws.Range("K7").Value = 7374.5
$ws.Range("L7").Value = 12397.6
$ws.Range("M7").Value = -7262.5
$ws.Range("N7").Value = -12621.6
$ws.Range("H16").Value = 47064.832
$ws.Range("I16").Value = 28830
$ws.Range("K16").Value = 28830
$ws.Range("M16").Value = -28660
$ws.Range("H22").Value = 3991
$ws.Range("I22").Value = 3991
$ws.Range("K22").Value = 3991
$ws.Range("M22").Value = -3696
$ws.Range("H27").Value = 3991
$ws.Range("I27").Value = 3991
$ws.Range("K27").Value = 3991
$ws.Range("M27").Value = -3884
$ws.Range("H46").Value = 1685.4839
$ws.Range("I46").Value = 1048.1538
$ws.Range("J46").Value = 2145.7778
$ws.Range("K46").Value = 1048.1538
$ws.Range("L46").Value = 2145.7778
$ws.Range("M46").Value = -860.1538
$ws.Range("N46").Value = -2521.7778
$ws.Range("H93").Value = 2748.6428
$ws.Range("I93").Value = 3168.7144
$ws.Range("J93").Value = 1488.4286
$ws.Range("K93").Value = 3168.7144
$ws.Range("L93").Value = 1488.4286
$ws.Range("M93").Value = -1920.7144
$ws.Range("N93").Value = -3984.4286
$ws.Range("H122").Value = 2677.2083
$ws.Range("I122").Value = 1846.7142
$ws.Range("J122").Value = 3839.9
$ws.Range("K122").Value = 5540.142599999999
$ws.Range("L122").Value = 11519.7
$ws.Range("M122").Value = -3090.142599999999
$ws.Range("N122").Value = -16419.7
$ws.Range("H126").Value = 9048.866
$ws.Range("I126").Value = 7374.5
$ws.Range("J126").Value = 12397.6
$ws.Range("K126").Value = 22123.5
$ws.Range("L126").Value = 37192.8
$ws.Range("M126").Value = -19653.5
$ws.Range("N126").Value = -42132.8

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H37").Value = 37500
$ws.Range("I37").Value = 37500
$ws.Range("K37").Value = 37500
$ws.Range("M37").Value = -37297
$ws.Range("H107").Value = 7409152
$ws.Range("I107").Value = 1352.1
$ws.Range("J107").Value = 22224752
$ws.Range("K107").Value = 4056.3
$ws.Range("L107").Value = 66674256
$ws.Range("M107").Value = -2136.3
$ws.Range("N107").Value = -66678096
$ws.Range("H113").Value = 444.76315
$ws.Range("I113").Value = 250.57143
$ws.Range("K113").Value = 751.71429
$ws.Range("M113").Value = 1418.28571
$ws.Range("H122").Value = 8109.3
$ws.Range("I122").Value = 4871.7617
$ws.Range("J122").Value = 15663.556
$ws.Range("K122").Value = 14615.2851
$ws.Range("L122").Value = 46990.66800000001
$ws.Range("M122").Value = -12165.2851
$ws.Range("N122").Value = -51890.66800000001
